# issue #5: stock data output to json file
#
# Adds a new "property_category" column (value "stock" for every row) to
# the 股票 (Stock) sheet, positioned right after the "total" column and
# before the "date" column. Also fixes a couple of stray-space typos in
# company names and normalises a malformed numeric-looking string.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# --- fix stray-space typos in company names -------------------------------
$ws.Range("B2").Value  = "世紀民生科技股份有限公司"
$ws.Range("B11").Value = "華泰商業銀行股份有限公司"

# --- normalise the malformed "total" value for 鍊德科技股份有限公司 --------
$ws.Range("G7").Value = ".2000000"

# --- insert the new "property_category" column (H), shifting            --
# --- date / legislator_name / legislator_id one column to the right     --
$ws.Columns.Item(8).Insert()

$ws.Range("H1").Value = "property_category"
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 8).Value = "stock"
}
